$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose literal text would be mangled by Excel's automatic
# text->number coercion (trailing/leading zeros, scientific notation)
# must be pre-formatted as Text so the literal string round-trips exactly.
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "27.253.92"
$ws.Cells.Item(2, 5).Value = "  -1.08%  "
$ws.Cells.Item(3, 4).Value = "1.785.96"
$ws.Cells.Item(3, 5).Value = "  +0.89%  "
$ws.Cells.Item(4, 4).Value = "1.007"
$ws.Cells.Item(4, 5).Value = "  +1.52%  "
$ws.Cells.Item(5, 4).Value = "336.48"
$ws.Cells.Item(5, 5).Value = "  +0.20%  "
$ws.Cells.Item(6, 4).Value = "1.006"
$ws.Cells.Item(6, 5).Value = "  +1.33%  "
$ws.Cells.Item(7, 4).Value = "0.3793"
$ws.Cells.Item(7, 5).Value = "  +0.52%  "
$ws.Cells.Item(8, 4).Value = "0.3430"
$ws.Cells.Item(8, 5).Value = "  -0.73%  "
$ws.Cells.Item(9, 4).Value = "48.45"
$ws.Cells.Item(9, 5).Value = "  -0.88%  "
$ws.Cells.Item(10, 4).Value = "1.193"
$ws.Cells.Item(10, 5).Value = "  -1.53%  "
$ws.Cells.Item(11, 4).Value = "0.07461"
$ws.Cells.Item(11, 5).Value = "  -2.52%  "
$ws.Cells.Item(12, 5).Value = "  +1.46%  "
$ws.Cells.Item(13, 4).Value = "21.80"
$ws.Cells.Item(13, 5).Value = "  +3.13%  "
$ws.Cells.Item(14, 4).Value = "6.437"
$ws.Cells.Item(14, 5).Value = "  -1.71%  "
$ws.Cells.Item(15, 4).Value = "1.785.04"
$ws.Cells.Item(15, 5).Value = "  +1.47%  "
$ws.Cells.Item(16, 4).Value = "7.067"
$ws.Cells.Item(16, 5).Value = "  -1.43%  "
$ws.Cells.Item(17, 4).Value = "0.00001098"
$ws.Cells.Item(17, 5).Value = "  -0.74%  "
$ws.Cells.Item(18, 4).Value = "0.06686"
$ws.Cells.Item(18, 5).Value = "  -0.92%  "
$ws.Cells.Item(19, 4).Value = "84.52"
$ws.Cells.Item(19, 5).Value = "  +0.19%  "
$ws.Cells.Item(20, 5).Value = "  +1.24%  "
$ws.Cells.Item(21, 4).Value = "6.528"
$ws.Cells.Item(21, 5).Value = "  +3.34%  "
$ws.Cells.Item(22, 4).Value = "17.35"
$ws.Cells.Item(22, 5).Value = "  +0.76%  "
$ws.Cells.Item(23, 4).Value = "27.227.69"
$ws.Cells.Item(23, 5).Value = "  -0.79%  "
$ws.Cells.Item(24, 5).Value = "  -4.49%  "
$ws.Cells.Item(25, 4).Value = "2.423"
$ws.Cells.Item(25, 5).Value = "  -1.12%  "
$ws.Cells.Item(26, 4).Value = "1.494"
$ws.Cells.Item(26, 5).Value = "  -0.47%  "
$ws.Cells.Item(27, 4).Value = "2.544"
$ws.Cells.Item(27, 5).Value = "  +2.60%  "
$ws.Cells.Item(28, 4).Value = "21.36"
$ws.Cells.Item(28, 5).Value = "  +6.73%  "
$ws.Cells.Item(29, 4).Value = "153.01"
$ws.Cells.Item(29, 5).Value = "  -0.03%  "
$ws.Cells.Item(30, 4).Value = "1.987.75"
$ws.Cells.Item(30, 5).Value = "  +1.54%  "
$ws.Cells.Item(31, 4).Value = "133.04"
$ws.Cells.Item(31, 5).Value = "  -0.84%  "
$ws.Cells.Item(32, 4).Value = "4.064"
$ws.Cells.Item(32, 5).Value = "  -1.32%  "
$ws.Cells.Item(33, 4).Value = "6.044"
$ws.Cells.Item(33, 5).Value = "  -2.24%  "
$ws.Cells.Item(34, 4).Value = "0.08646"
$ws.Cells.Item(34, 5).Value = "  -0.31%  "
$ws.Cells.Item(35, 5).Value = "  -0.96%  "
$ws.Cells.Item(36, 4).Value = "1.653"
$ws.Cells.Item(36, 5).Value = "  -2.87%  "
$ws.Cells.Item(37, 4).Value = "5.452"
$ws.Cells.Item(37, 5).Value = "  -1.49%  "
$ws.Cells.Item(38, 4).Value = "0.6855"
$ws.Cells.Item(38, 5).Value = "  +2.45%  "
$ws.Cells.Item(39, 4).Value = "0.06387"
$ws.Cells.Item(39, 5).Value = "  -0.54%  "
$ws.Cells.Item(40, 4).Value = "8.811"
$ws.Cells.Item(40, 5).Value = "  +1.54%  "
$ws.Cells.Item(41, 4).Value = "0.2192"
$ws.Cells.Item(41, 5).Value = "  -1.19%  "
$ws.Cells.Item(42, 4).Value = "0.02336"
$ws.Cells.Item(42, 5).Value = "  -2.18%  "
$ws.Cells.Item(43, 4).Value = "1.261"
$ws.Cells.Item(43, 5).Value = "  +2.28%  "
$ws.Cells.Item(44, 4).Value = "14.49"
$ws.Cells.Item(44, 5).Value = "  -0.52%  "
$ws.Cells.Item(45, 4).Value = "1.005"
$ws.Cells.Item(45, 5).Value = "  +1.29%  "
$ws.Cells.Item(46, 4).Value = "0.6411"
$ws.Cells.Item(46, 5).Value = "  +0.70%  "
$ws.Cells.Item(47, 4).Value = "3.858"
$ws.Cells.Item(47, 5).Value = "  -2.42%  "
$ws.Cells.Item(48, 4).Value = "2.128"
$ws.Cells.Item(48, 5).Value = "  +0.15%  "
$ws.Cells.Item(49, 4).Value = "128.98"
$ws.Cells.Item(49, 5).Value = "  -1.40%  "
$ws.Cells.Item(50, 4).Value = "0.07180"
$ws.Cells.Item(50, 5).Value = "  -1.75%  "
$ws.Cells.Item(51, 4).Value = "79.21"
$ws.Cells.Item(51, 5).Value = "  -0.20%  "

Write-Output "applied crypto price/volume updates"
